$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, copying the formatting used by the
# other header cells (e.g. G1: bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "Save"
$excel.CutCopyMode = 0

# Populate H2:H37 based on the "sum" column (G): 1 if sum > 10, else 0
for ($r = 2; $r -le 37; $r++) {
    $g = [double]$ws.Cells.Item($r, 7).Value2
    if ($g -gt 10) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
